$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# Step 1: Convert column D values in rows 25-35 from text to numeric.
for ($r = 25; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = [double]$cell.Value()
}

# Step 2: Append rows 36-46 as a copy of rows 25-35 (refreshed scrape),
# keeping column D as text and updating the Date Time column.
$srcStart = 25
$destStart = 36
$count = 11

for ($i = 0; $i -lt $count; $i++) {
    $srcRow = $srcStart + $i
    $destRow = $destStart + $i

    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value()
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value()
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value()

    # Column D must stay as text (not numeric) in the new rows.
    $dVal = [string]([int]$ws.Cells.Item($srcRow, 4).Value())
    $ws.Cells.Item($destRow, 4).Value = "'" + $dVal
    $ws.Cells.Item($destRow, 4).ClearFormats()

    $ws.Cells.Item($destRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value()
    $ws.Cells.Item($destRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value()
    $ws.Cells.Item($destRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value()
    $ws.Cells.Item($destRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value()
    $ws.Cells.Item($destRow, 9).Value = "17/06/2024 11:32:20"
}
